$d = $word.ActiveDocument

# --- locate the target paragraph -------------------------------------------------
# The paragraph that currently reads:
#   "(Ultimately Python provided ... coding frame" + "-" + <bookmark _GoBack> +
#   "work to carry out ... using python" + ")"
# needs to be "uncommented": drop the leading "(" / trailing ")", splice the
# hyphenated "frame-work" back together, add a trailing period, and split the
# paragraph in two - an (now empty) paragraph that keeps the _GoBack bookmark,
# followed by a fresh paragraph holding the restated sentence.

$target = $null
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ultimately Python provided a reasonable middle ground*") {
        $target = $p
        $targetIndex = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Ultimately Python provided...' paragraph"
}

$wasLastParagraph = ($target.Range.End -eq $d.Content.End)
$countBefore = $d.Paragraphs.Count

# --- build the replacement OOXML --------------------------------------------------
$newXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
    '<w:p>' +
      '<w:pPr><w:contextualSpacing w:val="0"/></w:pPr>' +
      '<w:bookmarkStart w:id="2" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="2"/>' +
    '</w:p>' +
    '<w:p>' +
      '<w:r><w:t>Ultimately Python provided a reasonable middle ground between team members knowledge and the potential difficulty in setting up the relevant infrastructure and coding frame-work to carry out the project. This middle ground led to the logical choice in deciding to imp</w:t></w:r>' +
      '<w:r><w:t>lement our project using python.</w:t></w:r>' +
    '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$null = $target.Range.InsertXML($newXml)

# InsertXML replaces the *content* of the range but - when the target paragraph is
# the very last paragraph of the document body - the original (now empty) trailing
# paragraph mark cannot be removed by a plain Range.Delete() on a zero-length range
# at end-of-story, so a stray empty paragraph is left behind after the two new ones.
# Detect that and fold it back in by deleting the paragraph mark that now separates
# the freshly-inserted final paragraph from the leftover empty one.
if ($wasLastParagraph -and $d.Paragraphs.Count -eq ($countBefore + 2)) {
    $lastReal = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $markRange = $d.Range($lastReal.Range.End - 1, $lastReal.Range.End)
    $null = $markRange.Delete()
}
